# Applies cryptos list update (Mon Jan  1 15:27:58 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.715.32"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "2.308.43"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.86"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.86"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.42%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.605"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.04"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0913"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.38"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.00%  "
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.991"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.33"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").Value = "2.651.65"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "2.297.98"
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("D18").Value = "42.843.42"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.44"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.83%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000105"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.52"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.53"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.49"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.53"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.88"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +19.25%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.93"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.30"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.12"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.30"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.74"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0869"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.85%  "
$ws.Range("E34").Value = "  +8.77%  "
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.64"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.33%  "
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("E39").Value = "  +2.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.64"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.80"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +10.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.57"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.03"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.230"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.29"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.03%  "
$ws.Range("D47").Value = "1.698.11"
$ws.Range("E47").Value = "  +1.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "111.10"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.32"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.86"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("E51").Value = "  -2.59%  "
